$d = $word.ActiveDocument

# Remove the first three paragraphs: two blank paragraphs followed by
# the "Instructions:" paragraph. We delete from the very start of the
# document through the end of the "Instructions:" paragraph (including
# its trailing paragraph mark), leaving the remaining content untouched.

$start = $d.Paragraphs(1).Range.Start
$end = $d.Paragraphs(3).Range.End

$r = $d.Range($start, $end)
$r.Delete()
